$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Purchase 22-23")
$ws2 = $wb.Worksheets.Item("Sale 22-23")

# Sheet1 ("Purchase 22-23"): delete rows 32-33 (the "Max International" entries)
$ws1.Rows("32:33").Delete()

# Sheet2 ("Sale 22-23"): delete rows 9-10 (old stale invoice entries)
$ws2.Rows("9:10").Delete()

# Update what is now row 10 (previously row 12) with new invoice data
$ws2.Range("B10").Value = 45173
$ws2.Range("C10").Value = "b23-24MQ207"
$ws2.Range("E10").Value = 57006
$ws2.Range("F10").Formula = "=E5+E6+E7+E8+E9+E10"

# Update selections / active sheet
$ws1.Range("C44").Select()
$ws2.Range("K10").Select()
$ws2.Activate()

Write-Output "done"
